# logstach conf and test data
# Updates the "data process.xlsx" benchmark sheet:
#  - L29 is overwritten with a literal space (was "=F29")
#  - the final summary row (33) gets refreshed benchmark numbers
#    (downstream SUM/weighted-average formulas in I33:M33 recalc automatically)
#  - the view scrolls down and the active selection moves to L39

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# L29: replace the old "=F29" formula with a single space value.
$ws.Range("L29").Value = " "

# Row 33 ("JDBC 0NF" benchmark row): refreshed raw measurements.
$ws.Range("C33").Value = 482683
$ws.Range("D33").Value = 42.74
$ws.Range("E33").Value = 22.51
$ws.Range("F33").Value = 19.15
$ws.Range("G33").Value = 2046

# Scroll the view down and move the selection to L39.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L39").Select()
